# Apply "nuevos experimentos no convexos" updates.
# Only the textual contents of a handful of (text-typed) cells change;
# cell types/styles/layout must stay exactly as before. Many of the new
# values look numeric (e.g. "-2.9"), so we temporarily force the Text
# number format while assigning them, then clear that formatting again
# so the cells keep behaving like plain shared-string cells.
# (Multi-area "A1,B2" ranges are avoided below -- they were observed to
# only apply NumberFormat/ClearFormats reliably to their first area.)

$wb = $excel.ActiveWorkbook

# NOTE: worksheet name lookup in this runtime is case-insensitive, and
# this workbook has both "Vector_bf" and "Vector_BF" sheets, which would
# collide if looked up by name. Use the (1-based) tab position instead,
# which matches the fixed sheet order of this workbook:
#   1 Funciones_Objetivo, 2 Restricciones_del_lider,
#   3 Restricciones_del_follower, 4 Punto_modificado,
#   5 Vector_bf, 6 Vector_BF, 7 Vector_Alpha
$wsLider    = $wb.Worksheets.Item(2)   # Restricciones_del_lider
$wsFollower = $wb.Worksheets.Item(3)   # Restricciones_del_follower
$wsPunto    = $wb.Worksheets.Item(4)   # Punto_modificado
$wsVecbf    = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBF    = $wb.Worksheets.Item(6)   # Vector_BF

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Restricciones_del_lider
Set-TextValue $wsLider "A2" "1.9 - x"
Set-TextValue $wsLider "B2" "-2.9"
Set-TextValue $wsLider "D2" "0.83"
Set-TextValue $wsLider "A3" "-1.9 + x"
Set-TextValue $wsLider "B3" "0.8999999999999999"
Set-TextValue $wsLider "D3" "0.08"

# Restricciones_del_follower
Set-TextValue $wsFollower "A2" "0.1499999999999999 - y"
Set-TextValue $wsFollower "B2" "-1.15"
Set-TextValue $wsFollower "D2" "0.01"
Set-TextValue $wsFollower "E2" "2.8000000000000003"
Set-TextValue $wsFollower "F2" "8.100000000000001"
Set-TextValue $wsFollower "A3" "-0.15000000000000002 + y"
Set-TextValue $wsFollower "B3" "-0.85"
Set-TextValue $wsFollower "D3" "0.97"
Set-TextValue $wsFollower "E3" "7.1"
Set-TextValue $wsFollower "F3" "2.1"

# Punto_modificado
Set-TextValue $wsPunto "A2" "1.9"
Set-TextValue $wsPunto "B2" "0.15"

# Vector_bf
Set-TextValue $wsVecbf "A2" "-0.39675000000000005"

# Vector_BF
Set-TextValue $wsVecBF "A2" "1.6"
Set-TextValue $wsVecBF "A3" "-9.2"
